{"js": "// Applies the edits described by the diff:\n//   1) \"tret\u00e5ig hackspett\" -> \"spillkr\u00e5ka och tret\u00e5ig hackspett\" (both occurrences)\n//   2) \"Detta \u00e4r en prioriterad art\" -> \"Dessa \u00e4r prioriterade arter\"\n//   3) \" denna art\" -> \" dessa arter\"\n//   4) \" arten\" -> \" arterna\"\n//   5) \"2026-02-10\" -> \"2026-02-11\"\n\nconst body = context.document.body;\n\n// 1) Replace every occurrence of \"tret\u00e5ig hackspett\" with\n//    \"spillkr\u00e5ka och tret\u00e5ig hackspett\" (there are two: one in the main\n//    paragraph, one further down in the bulleted recommendation).\nconst speciesRanges = body.search(\"tret\u00e5ig hackspett\", { matchCase: true });\nspeciesRanges.load(\"items\");\nawait context.sync();\n\nfor (const range of speciesRanges.items) {\n  range.insertText(\"spillkr\u00e5ka och tret\u00e5ig hackspett\", Word.InsertLocation.replace);\n}\n\n// 2) \"Detta \u00e4r en prioriterad art\" -> \"Dessa \u00e4r prioriterade arter\"\nconst singularStatement = body.search(\"Detta \u00e4r en prioriterad art\", { matchCase: true });\nsingularStatement.load(\"items\");\nawait context.sync();\n\nfor (const range of singularStatement.items) {\n  range.insertText(\"Dessa \u00e4r prioriterade arter\", Word.InsertLocation.replace);\n}\n\n// 3) \" denna art\" -> \" dessa arter\"\nconst thisSpecies = body.search(\" denna art\", { matchCase: true });\nthisSpecies.load(\"items\");\nawait context.sync();\n\nfor (const range of thisSpecies.items) {\n  range.insertText(\" dessa arter\", Word.InsertLocation.replace);\n}\n\n// 4) \" arten\" -> \" arterna\"\nconst theSpecies = body.search(\" arten\", { matchCase: true });\ntheSpecies.load(\"items\");\nawait context.sync();\n\nfor (const range of theSpecies.items) {\n  range.insertText(\" arterna\", Word.InsertLocation.replace);\n}\n\n// 5) Update the date at the bottom of the letter.\nconst dateRanges = body.search(\"2026-02-10\", { matchCase: true });\ndateRanges.load(\"items\");\nawait context.sync();\n\nfor (const range of dateRanges.items) {\n  range.insertText(\"2026-02-11\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Applies the edits described by the diff:\n#   1) \"tret\u00e5ig hackspett\" -> \"spillkr\u00e5ka och tret\u00e5ig hackspett\" (both occurrences)\n#   2) \"Detta \u00e4r en prioriterad art\" -> \"Dessa \u00e4r prioriterade arter\"\n#   3) \" denna art\" -> \" dessa arter\"\n#   4) \" arten\" -> \" arterna\"\n#   5) \"2026-02-10\" -> \"2026-02-11\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-AllText \"tret\u00e5ig hackspett\" \"spillkr\u00e5ka och tret\u00e5ig hackspett\"\nReplace-AllText \"Detta \u00e4r en prioriterad art\" \"Dessa \u00e4r prioriterade arter\"\nReplace-AllText \" denna art\" \" dessa arter\"\nReplace-AllText \" arten\" \" arterna\"\nReplace-AllText \"2026-02-10\" \"2026-02-11\"\n"}
